$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# Updated wind direction (wind_dir_fg) values
$wsFBS.Range("Q9").Value = "N"
$wsFBS.Range("Q24").Value = "E"
$wsFBS.Range("Q33").Value = "SW"
$wsFBS.Range("Q50").Value = "WNW"
$wsFBS.Range("Q52").Value = "SW"
$wsOther.Range("S12").Value = "WNW"

# Refresh the Timestamp column (AK) on the FBS sheet for every data row
$newTimestamp = "2024-10-13T12:04:29.414632"
for ($row = 2; $row -le 57; $row++) {
    $wsFBS.Range("AK$row").Value = $newTimestamp
}
